$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 7 (Tesla / Model 4 / 0 / 0) - entirely empty/placeholder entry
$ws.Rows.Item(7).Delete()

# Delete rows 3 (Nissan / Leaf) and 4 (Ford / Mustang Mach-E)
# so that Hyundai / Rivian shift up into rows 3 and 4
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(3).Delete()

$wb.Save()
